$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 27 (shifting the existing rows 27..31 down to 28..32)
$ws.Rows("27").Insert()

# Copy the cell formatting from row 28 (the product row just below) onto the
# newly inserted row 27 so it matches the look of the other product rows.
$ws.Range("A28:Q28").Copy()
$ws.Range("A27:Q27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new product row: "صابون ديتول اوريجنيال 115 جم"
$ws.Range("A27").Value = 21
$ws.Range("C27").Value = "صابون ديتول اوريجنيال 115 جم"
$ws.Range("H27").Value = "15:0"
$ws.Range("L27").Value = "0"
$ws.Range("N27").Value = "30.00"
$ws.Range("P27").Value = "30.0000"
$ws.Range("Q27").Value = "1:0"

# Re-create the merged cells for the new row 27 (same pattern as other rows)
$ws.Range("A27:B27").Merge()
$ws.Range("C27:G27").Merge()
$ws.Range("H27:K27").Merge()
$ws.Range("L27:M27").Merge()
$ws.Range("N27:O27").Merge()

# Update the grand total (was 726.68, now +30 for the new product = 756.68)
$ws.Range("P31").Value = 756.68

# Update the generated timestamp footer text
$ws.Range("A32").Value = "Wednesday, 20 August, 2025 2:15 PM"
